# Apply the Alvearie FHIR IG "documented-source" StructureDefinition refresh:
#  - bump Version 5.0.0 -> 6.0.0
#  - bump Date to the new publication timestamp
#  - fill in Publisher ("Alvearie Team")
#  - fix the duplicated "Contact" metadata row into "Jurisdiction"/"United States of America"
#  - drop the other stray duplicate "Contact" row entirely
#  - correct the root Extension element's Short/Definition text on the Elements sheet

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B3").Value = "6.0.0"
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"
$meta.Range("B9").Value = "Alvearie Team"
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Row 11 was a leftover duplicate of the old "Contact" row; remove it so the
# sheet collapses back down to 20 rows (A1:B20).
$meta.Rows.Item(11).Delete()

$elements = $wb.Worksheets.Item("Elements")
$elements.Range("K2").Value = "Documented Source"
$elements.Range("L2").Value = "Source type of the legal document"
